$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 29236.715
$ws.Range("J3").Value = 29236.715
$ws.Range("L3").Value = 29236.715
$ws.Range("N3").Value = -29464.715

$ws.Range("H12").Value = 220
$ws.Range("I12").Value = 166.5
$ws.Range("J12").Value = 300.25
$ws.Range("K12").Value = 166.5
$ws.Range("L12").Value = 300.25
$ws.Range("M12").Value = 3.5
$ws.Range("N12").Value = -640.25

$ws.Range("H40").Value = 3500
$ws.Range("J40").Value = 3500
$ws.Range("L40").Value = 3500
$ws.Range("N40").Value = -3850

$ws.Range("H63").Value = 31454
$ws.Range("J63").Value = 31454
$ws.Range("L63").Value = 31454
$ws.Range("N63").Value = -32702

$ws.Range("H64").Value = 3718.182
$ws.Range("I64").Value = 3342.8572
$ws.Range("J64").Value = 4375
$ws.Range("K64").Value = 3342.8572
$ws.Range("L64").Value = 4375
$ws.Range("M64").Value = -3094.8572
$ws.Range("N64").Value = -4871

$ws.Range("H66").Value = 31454
$ws.Range("J66").Value = 31454
$ws.Range("L66").Value = 94362
$ws.Range("N66").Value = -100602

$ws.Range("H67").Value = 3718.182
$ws.Range("I67").Value = 3342.8572
$ws.Range("J67").Value = 4375
$ws.Range("K67").Value = 3342.8572
$ws.Range("L67").Value = 4375
$ws.Range("M67").Value = -2484.8572
$ws.Range("N67").Value = -6091

$ws.Range("H74").Value = 3894.2942
$ws.Range("I74").Value = 3681.5
$ws.Range("J74").Value = 4198.2856
$ws.Range("K74").Value = 3681.5
$ws.Range("L74").Value = 4198.2856
$ws.Range("M74").Value = -2745.5
$ws.Range("N74").Value = -6070.2856

$ws.Range("H77").Value = 3894.2942
$ws.Range("I77").Value = 3681.5
$ws.Range("J77").Value = 4198.2856
$ws.Range("K77").Value = 18407.5
$ws.Range("L77").Value = 20991.428
$ws.Range("M77").Value = -13727.5
$ws.Range("N77").Value = -30351.428

$ws.Range("H86").Value = 16615
$ws.Range("I86").Value = 27490
$ws.Range("J86").Value = 5740
$ws.Range("K86").Value = 27490
$ws.Range("L86").Value = 5740
$ws.Range("M86").Value = -26367
$ws.Range("N86").Value = -7986

$ws.Range("H89").Value = 16615
$ws.Range("I89").Value = 27490
$ws.Range("J89").Value = 5740
$ws.Range("K89").Value = 137450
$ws.Range("L89").Value = 28700
$ws.Range("M89").Value = -131834
$ws.Range("N89").Value = -39932

$ws.Range("H92").Value = 516.05884
$ws.Range("I92").Value = 492.0625
$ws.Range("J92").Value = 900
$ws.Range("K92").Value = 492.0625
$ws.Range("L92").Value = 900
$ws.Range("M92").Value = 755.9375
$ws.Range("N92").Value = -3396

$ws.Range("H93").Value = 70694.31
$ws.Range("J93").Value = 70694.31
$ws.Range("L93").Value = 70694.31
$ws.Range("N93").Value = -75686.31

$ws.Range("H102").Value = 29236.715
$ws.Range("J102").Value = 29236.715
$ws.Range("L102").Value = 29236.715
$ws.Range("N102").Value = -35726.715

$ws.Range("H103").Value = 25000988
$ws.Range("I103").Value = 924.75
$ws.Range("J103").Value = 50001050
$ws.Range("K103").Value = 2774.25
$ws.Range("L103").Value = 150003150
$ws.Range("M103").Value = -2188.25
$ws.Range("N103").Value = -150004322

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 48200.715
$ws.Range("I107").Value = 59260.293
$ws.Range("J107").Value = 1197.5
$ws.Range("K107").Value = 59260.293
$ws.Range("L107").Value = 1197.5
$ws.Range("M107").Value = -57340.293
$ws.Range("N107").Value = -5037.5

$ws.Range("H108").Value = 39888
$ws.Range("J108").Value = 39888
$ws.Range("L108").Value = 39888
$ws.Range("N108").Value = -47568

$ws.Range("H110").Value = 44000
$ws.Range("J110").Value = 44000
$ws.Range("L110").Value = 44000
$ws.Range("N110").Value = -52180

$ws.Range("H112").Value = 2437.5
$ws.Range("I112").Value = 850
$ws.Range("J112").Value = 2543.3333
$ws.Range("K112").Value = 2550
$ws.Range("L112").Value = 7629.999899999999
$ws.Range("M112").Value = -1442
$ws.Range("N112").Value = -9845.999899999999

$ws.Range("H113").Value = 31253624
$ws.Range("I113").Value = 45456470
$ws.Range("J113").Value = 7357.6
$ws.Range("K113").Value = 45456470
$ws.Range("L113").Value = 7357.6
$ws.Range("M113").Value = -45453216
$ws.Range("N113").Value = -13865.6

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H118").Value = 50401380
$ws.Range("I118").Value = 84000240
$ws.Range("J118").Value = 3097
$ws.Range("K118").Value = 252000720
$ws.Range("L118").Value = 9291
$ws.Range("M118").Value = -251999063
$ws.Range("N118").Value = -12605

$ws.Range("H120").Value = 35000
$ws.Range("J120").Value = 35000
$ws.Range("L120").Value = 35000
$ws.Range("N120").Value = -44676

$ws.Range("H138").Value = 2540.1702
$ws.Range("I138").Value = 1393.5834
$ws.Range("J138").Value = 3736.6086
$ws.Range("K138").Value = 4180.7502
$ws.Range("L138").Value = 11209.8258
$ws.Range("M138").Value = 959.2497999999996
$ws.Range("N138").Value = -21489.8258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1829.2
$ws.Range("I122").Value = 1572.9333
$ws.Range("K122").Value = 4718.7999
$ws.Range("M122").Value = -2268.7999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 30000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 30000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 30000
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -30884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 18900
$ws.Range("J120").Value = 18900
$ws.Range("L120").Value = 56700
$ws.Range("N120").Value = -66376

$ws.Range("H122").Value = 2338.62
$ws.Range("I122").Value = 221
$ws.Range("J122").Value = 2935.8975
$ws.Range("K122").Value = 1989
$ws.Range("L122").Value = 26423.0775
$ws.Range("M122").Value = 461
$ws.Range("N122").Value = -31323.0775

$ws.Range("H123").Value = 4758.3335
$ws.Range("I123").Value = 1883.3334
$ws.Range("J123").Value = 5333.3335
$ws.Range("K123").Value = 5650.0002
$ws.Range("L123").Value = 16000.0005
$ws.Range("M123").Value = -3200.0002
$ws.Range("N123").Value = -20900.0005

$ws.Range("H125").Value = 2340.9092
$ws.Range("J125").Value = 2340.9092
$ws.Range("L125").Value = 7022.7276
$ws.Range("N125").Value = -16862.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1231.7894
$ws.Range("I46").Value = 875.125
$ws.Range("J46").Value = 3134
$ws.Range("K46").Value = 875.125
$ws.Range("L46").Value = 3134
$ws.Range("M46").Value = -687.125
$ws.Range("N46").Value = -3510

$ws.Range("H132").Value = 2886.1538
$ws.Range("I132").Value = 2052.2
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 6156.599999999999
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -3626.599999999999
$ws.Range("N132").Value = -22058
